# Update supersite names per Judi's spreadsheet 1-28-2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $val = $cell.Value2

    if ($val -eq "Trail Ridge MS") {
        $cell.Value2 = "Timberline K-8"
    }
    elseif ($val -eq "Manhattan MS") {
        $cell.Value2 = "Platt MS"
    }
    elseif ($val -eq "Monarch HS") {
        $cell.Value2 = "Monarch K8"
    }
    elseif ($val -eq "Ward Community Center") {
        $cell.Value2 = "Ward Town Hall"
    }
}
